$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, shifting the existing rows 39-90 down to 40-91.
$ws.Rows("39:39").Insert()

# Populate the newly inserted row 39 with the new data record.
$ws.Range("A39").Value2 = 8
$ws.Range("B39").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C39").Value2 = "Coquimbo"
$ws.Range("D39").Value2 = 44494
$ws.Range("E39").Value2 = 4
$ws.Range("F39").Value2 = 100112001
$ws.Range("G39").Value2 = "Berenjena"
$ws.Range("H39").Value2 = "Sin especificar"
$ws.Range("I39").Value2 = "Primera"
$ws.Range("J39").Value2 = 500
$ws.Range("K39").Value2 = 8000
$ws.Range("L39").Value2 = 9000
$ws.Range("M39").Value2 = 8500
$ws.Range("N39").Value2 = "`$/caja 60 unidades"
$ws.Range("O39").Value2 = "Región de Arica y Parinacota"
$ws.Range("P39").Value2 = 142
$ws.Range("Q39").Value2 = 60
$ws.Range("R39").Value2 = "Hortaliza"
